# Commit: "Added authorship, Fixed spacing in code"
# Appends a new account row (row 17) to the "accountInfo" sheet, matching
# the same unstyled inline-string layout used by the existing rows (e.g. row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("accountInfo")

$ws.Cells.Item(17, 1).Value = "random123"
$ws.Cells.Item(17, 2).Value = "Random1234!!"
$ws.Cells.Item(17, 3).Value = "narek.asaturyan@gmail.com"
$ws.Cells.Item(17, 4).Value = "The Witcher 3: Wild Hunt/Civilization VI/Half-Life: Alyx/Celeste/Stardew Valley/Doom"
